$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Reuniao - 1 hora" meeting entry that used to sit in D2
# (it is dropped from the weekly activities table; this also removes
# the now-unused shared string).
$ws.Range("D2").Value = ""

# Fix the accented character in the Friday row entry:
# "Reunião Sinavez - 1 hora" -> "Reuniao Sinavez - 1 hora"
$ws.Range("C6").Value = "Reuniao Sinavez - 1 hora"
